# Edit applied: refresh the 4 data rows with new sensor readings (commit: "custom accuracy + 데이터 1000개"),
# remove the now-obsolete last data row (old row 6), and nudge a handful of column widths.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 45143.50694444445
$ws.Range("B2").Value = 14.835
$ws.Range("C2").Value = 9.791
$ws.Range("D2").Value = 3.698
$ws.Range("E2").Value = 32.243
$ws.Range("F2").Value = 24.166
$ws.Range("G2").Value = 11.51
$ws.Range("H2").Value = 34.958
$ws.Range("I2").Value = 18.033
$ws.Range("J2").Value = 7.29
$ws.Range("K2").Value = 10.735
$ws.Range("L2").Value = 12.533
$ws.Range("M2").Value = 13.25
$ws.Range("N2").Value = 3.739
$ws.Range("O2").Value = 11.655
$ws.Range("P2").Value = 16.06
$ws.Range("Q2").Value = 10.282
$ws.Range("R2").Value = 3.096
$ws.Range("S2").Value = 1.74
$ws.Range("T2").Value = 170.025
$ws.Range("U2").Value = 32.298
$ws.Range("V2").Value = 10.758
$ws.Range("W2").Value = 20.812
$ws.Range("X2").Value = 10.713
$ws.Range("Y2").Value = 2.837
$ws.Range("Z2").Value = 18.288
$ws.Range("AA2").Value = 9.502
$ws.Range("AB2").Value = 8.642
$ws.Range("AC2").Value = 10.303
$ws.Range("AD2").Value = 12.679
$ws.Range("AE2").Value = 3.311
$ws.Range("AF2").Value = 31.418
$ws.Range("AG2").Value = 5.68
$ws.Range("AH2").Value = 13.449

# Row 3
$ws.Range("A3").Value = 45143.51388888889
$ws.Range("B3").Value = 4.751
$ws.Range("C3").Value = 2.968
$ws.Range("D3").Value = 1.387
$ws.Range("E3").Value = 10.618
$ws.Range("F3").Value = 7.407
$ws.Range("G3").Value = 3.622
$ws.Range("H3").Value = 17.485
$ws.Range("I3").Value = 5.817
$ws.Range("J3").Value = 2.304
$ws.Range("K3").Value = 3.1
$ws.Range("L3").Value = 4.085
$ws.Range("M3").Value = 4.358
$ws.Range("N3").Value = 1.216
$ws.Range("O3").Value = 3.76
$ws.Range("P3").Value = 5.15
$ws.Range("Q3").Value = 3.604
$ws.Range("R3").Value = 1.366
$ws.Range("S3").Value = 0.68
$ws.Range("T3").Value = 49.94
$ws.Range("U3").Value = 10.749
$ws.Range("V3").Value = 3.47
$ws.Range("W3").Value = 6.738
$ws.Range("X3").Value = 3.408
$ws.Range("Y3").Value = 1.118
$ws.Range("Z3").Value = 8.406
$ws.Range("AA3").Value = 3.065
$ws.Range("AB3").Value = 2.934
$ws.Range("AC3").Value = 3.473
$ws.Range("AD3").Value = 4.003
$ws.Range("AE3").Value = 1.246
$ws.Range("AF3").Value = 16.597
$ws.Range("AG3").Value = 1.709
$ws.Range("AH3").Value = 4.341

# Row 4
$ws.Range("A4").Value = 45143.52083333334
$ws.Range("B4").Value = 3.794
$ws.Range("C4").Value = 2.481
$ws.Range("D4").Value = 0.887
$ws.Range("E4").Value = 8.512
$ws.Range("F4").Value = 6.084
$ws.Range("G4").Value = 2.896
$ws.Range("H4").Value = 12.977
$ws.Range("I4").Value = 4.654
$ws.Range("J4").Value = 1.87
$ws.Range("K4").Value = 2.553
$ws.Range("L4").Value = 3.312
$ws.Range("M4").Value = 3.554
$ws.Range("N4").Value = 0.968
$ws.Range("O4").Value = 3.008
$ws.Range("P4").Value = 4.112
$ws.Range("Q4").Value = 2.856
$ws.Range("R4").Value = 0.898
$ws.Range("S4").Value = 0.44
$ws.Range("T4").Value = 38.449
$ws.Range("U4").Value = 8.469
$ws.Range("V4").Value = 2.776
$ws.Range("W4").Value = 5.336
$ws.Range("X4").Value = 2.763
$ws.Range("Y4").Value = 0.813
$ws.Range("Z4").Value = 6.084
$ws.Range("AA4").Value = 2.452
$ws.Range("AB4").Value = 2.333
$ws.Range("AC4").Value = 2.745
$ws.Range("AD4").Value = 3.274
$ws.Range("AE4").Value = 0.766
$ws.Range("AF4").Value = 12.048
$ws.Range("AG4").Value = 1.391
$ws.Range("AH4").Value = 3.471

# Row 5
$ws.Range("A5").Value = 45143.52777777778
$ws.Range("B5").Value = 18.69
$ws.Range("C5").Value = 13.78
$ws.Range("D5").Value = 1.18
$ws.Range("E5").Value = 40.89
$ws.Range("F5").Value = 33.07
$ws.Range("G5").Value = 14.64
$ws.Range("H5").Value = 52.38
$ws.Range("I5").Value = 22.69
$ws.Range("J5").Value = 10.01
$ws.Range("K5").Value = 14.78
$ws.Range("L5").Value = 16.33
$ws.Range("M5").Value = 17.34
$ws.Range("N5").Value = 4.71
$ws.Range("O5").Value = 14.66
$ws.Range("P5").Value = 20.8
$ws.Range("Q5").Value = 12.43
$ws.Range("R5").Value = 0.8
$ws.Range("S5").Value = 0.78
$ws.Range("T5").Value = 215.84
$ws.Range("U5").Value = 40.88
$ws.Range("V5").Value = 13.53
$ws.Range("W5").Value = 27.39
$ws.Range("X5").Value = 14.39
$ws.Range("Y5").Value = 2.24
$ws.Range("Z5").Value = 26.16
$ws.Range("AA5").Value = 11.95
$ws.Range("AB5").Value = 10.63
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 17.01
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 47.48
$ws.Range("AG5").Value = 7.57
$ws.Range("AH5").Value = 16.92

# Remove the old trailing row (previously row 6)
$ws.Rows.Item(6).Delete()

# Column width adjustments
$ws.Range("G:G").ColumnWidth = 6.15
$ws.Range("L:L").ColumnWidth = 7.15
$ws.Range("M:M").ColumnWidth = 6.15
$ws.Range("P:P").ColumnWidth = 6.15
$ws.Range("Q:Q").ColumnWidth = 7.15
$ws.Range("S:S").ColumnWidth = 5.15
$ws.Range("T:T").ColumnWidth = 8.15
$ws.Range("V:V").ColumnWidth = 7.15
$ws.Range("AC:AC").ColumnWidth = 7.15
